$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.713.96"
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = "'1.875.39"
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -1.14%  '
$ws.Range('D5').Value = "'314.63"
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('D7').Value = "'0.5084"
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('D8').Value = "'0.3907"
$ws.Range('E8').Value = '  -1.35%  '
$ws.Range('D9').Value = "'0.08367"
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('D10').Value = "'42.24"
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').Value = "'1.108"
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = "'6.188"
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').Value = "'1.873.25"
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = "'20.37"
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = "'7.253"
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = "'1.007"
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = "'93.16"
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('D18').Value = "'0.00001099"
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = "'0.06709"
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = "'17.64"
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').Value = "'5.933"
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = "'28.718.56"
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('D24').Value = "'11.10"
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('D25').Value = "'2.227"
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').Value = "'2.084.92"
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').Value = "'157.43"
$ws.Range('E27').Value = '  -2.76%  '
$ws.Range('D28').Value = "'20.70"
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = "'2.412"
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('D30').Value = "'126.33"
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('D31').Value = "'0.1037"
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('D32').Value = "'1.046"
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').Value = "'5.788"
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = "'3.644"
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').Value = "'0.02448"
$ws.Range('D36').Value = "'0.06532"
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('D37').Value = "'8.985"
$ws.Range('E37').Value = '  +1.19%  '
$ws.Range('D38').Value = "'0.2164"
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').Value = "'5.064"
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('D40').Value = "'1.194"
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').Value = "'1.238"
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').Value = "'0.6382"
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = "'11.17"
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('D45').Value = "'0.5991"
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('D46').Value = "'13.06"
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').Value = "'3.674"
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = "'2.005"
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = "'1.224"
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('D50').Value = "'122.06"
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').Value = '  -2.19%  '
